$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Coding Phase Defects" -- SONARLINT fix: System.out -> ApplicationLogger
# ---------------------------------------------------------------------------
$wsCoding = $wb.Worksheets.Item("Coding Phase Defects")
$wsCoding.Range("E11").Value = "changed System.out to ApplicationLogger.log"
$wsCoding.Range("E12").Value = "changed System.out to ApplicationLogger.log"

# ---------------------------------------------------------------------------
# Sheet: "DynamicCodeAnalysis" -- new SonarLint findings (rows 10-12)
# ---------------------------------------------------------------------------
$wsDyn = $wb.Worksheets.Item("DynamicCodeAnalysis")

# Row 10: Main.java/22 (typed Issue before File,Line)
$wsDyn.Range("D10").Value = "Standard outputs should not be used directly to log anything"
$wsDyn.Range("C10").Value = "Main.java/22"
$wsDyn.Range("E10").Value = "System.err.println(service.getAllParts());"
$wsDyn.Range("F10").Value = "ApplicationLogger.log(Level.INFO, service.getAllProducts().toString());"

# Row 11: InventoryService.java/47
$wsDyn.Range("C11").Value = "InventoryService.java/47"
$wsDyn.Range("D11").Value = "Methods should not have too many parameters"
$wsDyn.Range("E11").Value = "public void updateInhousePart(int partIndex, int partId, String name, double price, int inStock, int min, int max, int partDynamicValue)"
$wsDyn.Range("F11").Value = "changing ,method requires creating new classes to wrap all the parameters into one which is too costly timewise and requires redesign of all the UML diagrams"

# Row 12: AddProductController/28 (typed Issue/Before before File,Line)
$wsDyn.Range("D12").Value = 'Unused "private" fields should be removed'
$wsDyn.Range("E12").Value = "private int productId;"
$wsDyn.Range("C12").Value = "AddProductController/28"

# Row heights (rows 10 & 11 were manually resized to fit the new text)
$wsDyn.Rows.Item(10).RowHeight = 15
$wsDyn.Rows.Item(11).RowHeight = 30

# Column widths widened to fit the new SonarLint text
$wsDyn.Columns.Item(4).ColumnWidth = 39.736979166666664
$wsDyn.Columns.Item(5).ColumnWidth = 67.89760416666667
$wsDyn.Columns.Item(6).ColumnWidth = 73.89760416666667

# ---------------------------------------------------------------------------
# Row height on Coding sheet (auto-expanded because of the longer text)
# ---------------------------------------------------------------------------
$wsCoding.Rows.Item(11).RowHeight = 30
$wsCoding.Rows.Item(12).RowHeight = 30

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping: work finished on Coding sheet (H13)
# then moved to DynamicCodeAnalysis (C13), which ends up the active tab.
# ---------------------------------------------------------------------------
$wsCoding.Range("H13").Select()
$wsDyn.Activate()
$wsDyn.Range("C13").Select()

Write-Host "done"
